# Applies the "Conflict Fix and Penalty scaling to 0 to 100%" update to the
# TabuSearch_Stats workbook: refreshed per-division optimization timings
# (table 1), rescaled penalty/quality numbers (table 2) and the updated
# "Entire League" totals row (table 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Table 1 (rows 2-12): Attempted Moves / ... / Optimization Time ----
$table1 = @{
    2  = @{ C = 11049.0; E = 11049.0; F = "24 sec" }
    3  = @{ C = 838.0;   E = 838.0;   F = "6 sec" }
    4  = @{ C = 1036.0;  E = 1036.0;  F = "6 sec" }
    5  = @{ C = 5740.0;  E = 5736.0;  F = "20 sec" }
    6  = @{ C = 12046.0; E = 12046.0; F = "15 sec" }
    7  = @{ C = 15148.0; E = 15146.0; F = "17 sec" }
    8  = @{ C = 11196.0; E = 11196.0; F = "21 sec" }
    9  = @{ C = 884.0;   E = 884.0;   F = "5 sec" }
    10 = @{ C = 5055.0;  E = 5055.0;  F = "17 sec" }
    11 = @{ C = 1344.0;  E = 1344.0;  F = "8 sec" }
    12 = @{ C = 12495.0; E = 12495.0; F = "18 sec" }
}

foreach ($row in $table1.Keys) {
    $vals = $table1[$row]
    $ws.Range("C$row").Value2 = $vals.C
    $ws.Range("E$row").Value2 = $vals.E
    $ws.Range("F$row").Value2 = $vals.F
}

# ---- Table 2 (rows 16-26): rescaled penalty (Quality) numbers ----
$table2 = @{
    16 = @{ B = 92.55;  C = 1471.0; D = 9.0;  E = 1462.0 }
    17 = @{ B = 120.18; C = 947.0;  D = 11.0; E = 936.0 }
    18 = @{ B = 101.91; C = 3149.0; D = 10.0; E = 3139.0 }
    19 = @{ C = 900.0;  D = 7.0;    E = 893.0 }
    20 = @{ B = 157.0;  C = 456.0;  D = 22.0; E = 434.0 }
    21 = @{ B = 93.18;  C = 1298.0; D = 21.0; E = 1277.0 }
    22 = @{ C = 1242.0; D = 16.0;   E = 1226.0 }
    23 = @{ B = 107.82; C = 825.0;  D = 3.0;  E = 822.0 }
    24 = @{ B = 100.0;  C = 1224.0; D = 22.0; E = 1202.0 }
    25 = @{ B = 113.36; C = 1362.0; D = 13.0; E = 1349.0 }
    26 = @{ B = 209.18; C = 3839.0; D = 41.0; E = 3798.0 }
}

foreach ($row in $table2.Keys) {
    $vals = $table2[$row]
    if ($vals.ContainsKey("B")) { $ws.Range("B$row").Value2 = $vals.B }
    $ws.Range("C$row").Value2 = $vals.C
    $ws.Range("D$row").Value2 = $vals.D
    $ws.Range("E$row").Value2 = $vals.E
}

# ---- Table 3 (row 30): Entire League totals ----
$ws.Range("B30").Value2 = 359.17
$ws.Range("C30").Value2 = 16713.0
$ws.Range("D30").Value2 = 175.0
$ws.Range("E30").Value2 = 16538.0
$ws.Range("F30").Value2 = "2 min, 59 sec"
